$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.494.27"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.394.24"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'576.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'141.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D9").Value = "'7.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "3.975.42"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D15").Value = "3.401.77"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "61.481.93"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "'6.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'13.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").Value = "'8.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'391.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").Value = "'75.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.555"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'0.0000114"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Value = "'8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'2.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -4.06%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").Value = "'167.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'5.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").Value = "3.427.35"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "'0.0771"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'26.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.63%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "2.470.56"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "'22.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "'6.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "'2.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("E51").Value = "  -1.28%  "
